$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.808.57"
$ws.Range("E2").Value = "  +8.82%  "

$ws.Range("D3").Value = "'1.952.51"
$ws.Range("E3").Value = "  +7.39%  "

$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "'342.21"
$ws.Range("E5").Value = "  +2.97%  "

$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").Value = "'0.4781"
$ws.Range("E7").Value = "  +4.95%  "

$ws.Range("E8").Value = "  +9.09%  "

$ws.Range("D9").Value = "'48.51"
$ws.Range("E9").Value = "  +5.78%  "

$ws.Range("D10").Value = "'0.08269"
$ws.Range("E10").Value = "  +5.71%  "

$ws.Range("E11").Value = "  +9.14%  "

$ws.Range("D12").Value = "'22.70"
$ws.Range("E12").Value = "  +8.61%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.949.33"
$ws.Range("E13").Value = "  +8.02%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.211"
$ws.Range("E14").Value = "  +7.08%  "

$ws.Range("D15").Value = "'7.437"
$ws.Range("E15").Value = "  +5.92%  "

$ws.Range("D16").Value = "'92.45"
$ws.Range("E16").Value = "  +3.77%  "

$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("E18").Value = "  +4.89%  "

$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("E20").Value = "  +6.65%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").Value = "'29.762.86"
$ws.Range("E22").Value = "  +8.69%  "

$ws.Range("D23").Value = "'5.608"
$ws.Range("E23").Value = "  +6.67%  "

$ws.Range("E24").Value = "  +4.90%  "

$ws.Range("D25").Value = "'2.282"
$ws.Range("E25").Value = "  +1.17%  "

$ws.Range("D26").Value = "'2.171.59"
$ws.Range("E26").Value = "  +7.27%  "

$ws.Range("D27").Value = "'160.99"
$ws.Range("E27").Value = "  +1.82%  "

$ws.Range("D28").Value = "'20.21"
$ws.Range("E28").Value = "  +5.22%  "

$ws.Range("D29").Value = "'2.191"
$ws.Range("E29").Value = "  +8.01%  "

$ws.Range("E30").Value = "  +8.85%  "

$ws.Range("D31").Value = "'122.65"
$ws.Range("E31").Value = "  +4.45%  "

$ws.Range("D32").Value = "'1.019"
$ws.Range("E32").Value = "  +9.96%  "

$ws.Range("D33").Value = "'0.09638"
$ws.Range("E33").Value = "  +3.58%  "

$ws.Range("D34").Value = "'1.483"
$ws.Range("E34").Value = "  +13.54%  "

$ws.Range("D35").Value = "'3.682"
$ws.Range("E35").Value = "  +3.38%  "

$ws.Range("E36").Value = "  +6.44%  "

$ws.Range("D37").Value = "'0.06313"
$ws.Range("E37").Value = "  +7.37%  "

$ws.Range("D38").Value = "'0.02335"
$ws.Range("E38").Value = "  +7.23%  "

$ws.Range("D39").Value = "'8.549"
$ws.Range("E39").Value = "  +6.28%  "

$ws.Range("D40").Value = "'1.199"
$ws.Range("E40").Value = "  +5.42%  "

$ws.Range("D41").Value = "'0.6120"
$ws.Range("E41").Value = "  +7.26%  "

$ws.Range("E42").Value = "  +9.01%  "

$ws.Range("E43").Value = "  +5.24%  "

$ws.Range("D44").Value = "'0.9993"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").Value = "'1.292"
$ws.Range("E45").Value = "  +2.46%  "

$ws.Range("D46").Value = "'2.404"
$ws.Range("E46").Value = "  +33.43%  "

$ws.Range("D47").Value = "'12.62"
$ws.Range("E47").Value = "  +7.82%  "

$ws.Range("D48").Value = "'0.5724"
$ws.Range("E48").Value = "  +6.67%  "

$ws.Range("D49").Value = "'2.001"
$ws.Range("E49").Value = "  +7.46%  "

$ws.Range("D50").Value = "'0.07396"
$ws.Range("E50").Value = "  +12.73%  "

$ws.Range("D51").Value = "'114.51"
$ws.Range("E51").Value = "  +4.29%  "
